# Manual evaluation COP22, COP22
# Row 27 in Sheet1 is reverted back to an "unevaluated" row (like most other
# rows): Relevance -> "no", and the Topic/Unit/Scale/Time/Principle/30-word
# columns are cleared out. The row height shrinks back down to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relevance column: "yes" -> "no"
$ws.Range("B27").Value = "no"

# Topic (C27) and Unit (D27) are fully cleared (no leftover style either)
$ws.Range("C27:D27").Clear()

# Scale (E27), Time (F27), Principle (G27) and 30-word explanation (H27) are
# cleared but keep their existing (wrap-text) cell style
$ws.Range("E27:H27").ClearContents()

# Row shrinks from the old 150pt to 75pt now that it holds far less text
$ws.Rows.Item(27).RowHeight = 75

# Move the active selection to the cell that was most recently being worked on
$ws.Range("H27").Select()
